$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.872.79'
$ws.Range("E2").Value = '  +3.54%  '

$ws.Range("D3").Value = '2.264.27'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''299.56'
$ws.Range("E5").Value = '  -0.42%  '

$ws.Range("D6").Value = '''100.91'
$ws.Range("E6").Value = '  +6.00%  '

$ws.Range("D7").Value = '''0.559'
$ws.Range("E7").Value = '  -1.16%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '''0.509'
$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("D10").Value = '''35.39'
$ws.Range("E10").Value = '  +3.44%  '

$ws.Range("D11").Value = '''0.0784'
$ws.Range("E11").Value = '  -1.15%  '

$ws.Range("D12").Value = '''7.07'
$ws.Range("E12").Value = '  -1.83%  '

$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("D14").Value = '2.608.12'
$ws.Range("E14").Value = '  -0.13%  '

$ws.Range("D15").Value = '2.262.42'
$ws.Range("E15").Value = '  -0.10%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''13.58'
$ws.Range("E16").Value = '  -0.12%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '46.841.27'
$ws.Range("E17").Value = '  +3.93%  '

$ws.Range("D18").Value = '''0.791'
$ws.Range("E18").Value = '  -2.23%  '

$ws.Range("D19").Value = '''12.73'
$ws.Range("E19").Value = '  -4.18%  '

$ws.Range("D20").Value = '0.0₃0932'
$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("D21").Value = '''5.80'
$ws.Range("E21").Value = '  -3.53%  '

$ws.Range("D22").Value = '''65.22'
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").Value = '''246.14'
$ws.Range("E23").Value = '  +2.65%  '

$ws.Range("D24").Value = '''2.80'
$ws.Range("E24").Value = '  -3.08%  '

$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '''1.86'
$ws.Range("E26").Value = '  -2.44%  '

$ws.Range("D27").Value = '''41.77'
$ws.Range("E27").Value = '  +0.70%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").Value = '''9.65'
$ws.Range("E29").Value = '  +0.98%  '

$ws.Range("D30").Value = '''20.40'
$ws.Range("E30").Value = '  +4.05%  '

$ws.Range("D31").Value = '''2.83'
$ws.Range("E31").Value = '  +9.70%  '

$ws.Range("D32").Value = '''145.68'
$ws.Range("E32").Value = '  -4.63%  '

$ws.Range("D33").Value = '''3.28'
$ws.Range("E33").Value = '  +12.52%  '

$ws.Range("D34").Value = '''5.37'
$ws.Range("E34").Value = '  -2.48%  '

$ws.Range("D35").Value = '''0.0766'
$ws.Range("E35").Value = '  -2.80%  '

$ws.Range("E36").Value = '  +10.62%  '

$ws.Range("E37").Value = '  -2.44%  '

$ws.Range("D38").Value = '''15.97'
$ws.Range("E38").Value = '  +18.60%  '

$ws.Range("D39").Value = '''1.69'
$ws.Range("E39").Value = '  -4.07%  '

$ws.Range("D40").Value = '''3.85'
$ws.Range("E40").Value = '  -1.24%  '

$ws.Range("D41").Value = '''0.0296'
$ws.Range("E41").Value = '  -5.19%  '

$ws.Range("D42").Value = '''3.12'
$ws.Range("E42").Value = '  -2.85%  '

$ws.Range("D43").Value = '''0.998'
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").Value = '''91.24'
$ws.Range("E44").Value = '  +19.83%  '

$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '''1.90'
$ws.Range("E45").Value = '  -2.24%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.779.34'
$ws.Range("E46").Value = '  +0.80%  '

$ws.Range("D47").Value = '''71.29'
$ws.Range("E47").Value = '  +1.71%  '

$ws.Range("D48").Value = '''0.185'
$ws.Range("E48").Value = '  -3.86%  '

$ws.Range("D49").Value = '''4.81'
$ws.Range("E49").Value = '  +1.86%  '

$ws.Range("D50").Value = '2.486.92'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("D51").Value = '''7.80'
$ws.Range("E51").Value = '  -0.74%  '
